# Generate Report for Handback
# Refresh the handoff/handback timestamp strings recorded on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-09-06 05:21:43"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first file.
$wsZhCn.Range("H2").Value = "2016-09-06 05:21:38"
$wsZhCn.Range("K2").Value = "2016-09-06 05:22:11"

# de-de sheet: "Correspond Handoff Datetime" (shares the same timestamp text as the
# Overview sheet's generate date, so it moves together) and "Correspond Handback DateTime".
$wsDeDe.Range("H2").Value = "2016-09-06 05:21:43"
$wsDeDe.Range("K2").Value = "2016-09-06 05:22:20"
